$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283, shifting rows 283:320 down to 284:321.
$ws.Rows(283).Insert()

# Populate the newly inserted row 283 with the new weekly record.
# (Same market/region/category/variety/quality/unit/origin/classification
# as the record that used to sit at row 283, but with a new date, volume,
# price range and $/Kg figure.)
$ws.Cells.Item(283, 1).Value = 10
$ws.Cells.Item(283, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(283, 3).Value = "La Araucanía"
$ws.Cells.Item(283, 4).Value = 45131
$ws.Cells.Item(283, 5).Value = 9
$ws.Cells.Item(283, 6).Value = 100112005
$ws.Cells.Item(283, 7).Value = "Puerro"
$ws.Cells.Item(283, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 110
$ws.Cells.Item(283, 11).Value = 8000
$ws.Cells.Item(283, 12).Value = 8000
$ws.Cells.Item(283, 13).Value = 8000
$ws.Cells.Item(283, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(283, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(283, 16).Value = 667
$ws.Cells.Item(283, 17).Value = 12
$ws.Cells.Item(283, 18).Value = "Hortaliza"
